$d = $word.ActiveDocument

# 1. Remove the "Phan cong.PNG" picture from the first paragraph, leaving an
#    empty paragraph behind (the <w:p> stays, its <w:r><w:drawing>... run goes).
$d.InlineShapes.Item(1).Delete()

# 2. Update the "Sauk hi viet xong ..." paragraph text to the new note.
$p = $d.Paragraphs.Item(6)
$r = $p.Range
$r.MoveEnd(1, -1)
$r.Text = "Phúc: CHỉnh SRS"

# 3. Move the "_GoBack" bookmark from the last paragraph onto the paragraph
#    we just edited (collapsed, right after the new text).
$d.Bookmarks.Item("_GoBack").Delete()

$p = $d.Paragraphs.Item(6)
$lastChar = $d.Range($p.Range.End - 2, $p.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $lastChar)

# Collapse the newly-added bookmark down to a zero-length mark positioned
# right after the text (Bookmarks.Add with a pre-collapsed range lands at
# offset 0 in this runtime, so build it from a 1-char range and shrink it).
$bmRange = $d.Bookmarks.Item("_GoBack").Range
$savedText = $bmRange.Text
$bmRange.Text = ""
$d.Bookmarks.Item("_GoBack").Range.InsertBefore($savedText)
